$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 7

$ws.Range("B4").Select()
